# Scheduled market-data refresh for Seraph_Profits workbook.
# Updates currentAveragePrice* / Leve*Price* / Leve*Profit* columns (H:N)
# per crafting-profession sheet, reflecting latest Universalis price pulls.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 138.375
$ws.Range("J9").Value = 214.33333
$ws.Range("L9").Value = 214.33333
$ws.Range("N9").Value = -552.3333299999999
$ws.Range("H43").Value = 3653
$ws.Range("I43").Value = 2692
$ws.Range("J43").Value = 5094.5
$ws.Range("K43").Value = 2692
$ws.Range("L43").Value = 5094.5
$ws.Range("M43").Value = -2623
$ws.Range("N43").Value = -5232.5
$ws.Range("H55").Value = 630.3333
$ws.Range("I55").Value = 225.5
$ws.Range("J55").Value = 746
$ws.Range("K55").Value = 225.5
$ws.Range("L55").Value = 746
$ws.Range("M55").Value = -11.5
$ws.Range("N55").Value = -1174
$ws.Range("H58").Value = 2377.8333
$ws.Range("J58").Value = 3339
$ws.Range("L58").Value = 10017
$ws.Range("N58").Value = -10317
$ws.Range("H111").Value = 824.8
$ws.Range("J111").Value = 1466
$ws.Range("L111").Value = 4398
$ws.Range("N111").Value = -10532
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 1513.9546
$ws.Range("I132").Value = 1384.1666
$ws.Range("K132").Value = 4152.4998
$ws.Range("M132").Value = -1622.4998
$ws.Range("H135").Value = 1015.2105
$ws.Range("I135").Value = 528.7143
$ws.Range("K135").Value = 4758.428699999999
$ws.Range("M135").Value = -2223.428699999999
$ws.Range("H137").Value = 2317.0417
$ws.Range("I137").Value = 1229.7858
$ws.Range("J137").Value = 3839.2
$ws.Range("K137").Value = 3689.3574
$ws.Range("L137").Value = 11517.6
$ws.Range("M137").Value = -1139.3574
$ws.Range("N137").Value = -16617.6
$ws.Range("H138").Value = 3770.9253
$ws.Range("I138").Value = 1120.8
$ws.Range("K138").Value = 3362.4
$ws.Range("M138").Value = 1777.6

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5026.5835
$ws.Range("I74").Value = 1348.8334
$ws.Range("K74").Value = 1348.8334
$ws.Range("M74").Value = -474.8334
$ws.Range("H77").Value = 5026.5835
$ws.Range("I77").Value = 1348.8334
$ws.Range("K77").Value = 6744.166999999999
$ws.Range("M77").Value = -2376.166999999999
$ws.Range("H102").Value = 1135.1
$ws.Range("I102").Value = 821.4666999999999
$ws.Range("J102").Value = 2076
$ws.Range("K102").Value = 821.4666999999999
$ws.Range("L102").Value = 2076
$ws.Range("M102").Value = 800.5333000000001
$ws.Range("N102").Value = -5320
$ws.Range("H132").Value = 1099.6
$ws.Range("I132").Value = 1099.6
$ws.Range("K132").Value = 3298.8
$ws.Range("M132").Value = -768.7999999999997

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3588.5881
$ws.Range("I134").Value = 3286.5
$ws.Range("K134").Value = 9859.5
$ws.Range("M134").Value = -7324.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2000430.2
$ws.Range("I6").Value = 683.6667
$ws.Range("J6").Value = 5000050
$ws.Range("K6").Value = 683.6667
$ws.Range("L6").Value = 5000050
$ws.Range("M6").Value = -570.6667
$ws.Range("N6").Value = -5000276
$ws.Range("H48").Value = 9999
$ws.Range("J48").Value = 9999
$ws.Range("L48").Value = 9999
$ws.Range("N48").Value = -10951
$ws.Range("H58").Value = 3147.0476
$ws.Range("I58").Value = 1902.1111
$ws.Range("J58").Value = 4080.75
$ws.Range("K58").Value = 1902.1111
$ws.Range("L58").Value = 4080.75
$ws.Range("M58").Value = -1699.1111
$ws.Range("N58").Value = -4486.75
$ws.Range("H132").Value = 1866.2858
$ws.Range("I132").Value = 1493.8125
$ws.Range("J132").Value = 2362.9167
$ws.Range("K132").Value = 4481.4375
$ws.Range("L132").Value = 7088.750100000001
$ws.Range("M132").Value = -1951.4375
$ws.Range("N132").Value = -12148.7501
$ws.Range("H134").Value = 3881.5557
$ws.Range("I134").Value = 3331.4614
$ws.Range("K134").Value = 9994.3842
$ws.Range("M134").Value = -7459.3842
$ws.Range("H136").Value = 3147.0476
$ws.Range("I136").Value = 1902.1111
$ws.Range("J136").Value = 4080.75
$ws.Range("K136").Value = 5706.3333
$ws.Range("L136").Value = 12242.25
$ws.Range("M136").Value = -3156.3333
$ws.Range("N136").Value = -17342.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1322.9166
$ws.Range("I14").Value = 1322.9166
$ws.Range("K14").Value = 3968.7498
$ws.Range("M14").Value = -3795.7498
$ws.Range("H38").Value = 117.809525
$ws.Range("I38").Value = 33.294117
$ws.Range("J38").Value = 477
$ws.Range("K38").Value = 99.882351
$ws.Range("L38").Value = 1431
$ws.Range("M38").Value = 247.117649
$ws.Range("N38").Value = -2125
$ws.Range("H39").Value = 1822.6
$ws.Range("J39").Value = 1799.8
$ws.Range("L39").Value = 5399.4
$ws.Range("N39").Value = -5987.4
$ws.Range("H55").Value = 85358.25
$ws.Range("J55").Value = 2862.375
$ws.Range("L55").Value = 8587.125
$ws.Range("N55").Value = -8941.125

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2811.3333
$ws.Range("I80").Value = 2112
$ws.Range("K80").Value = 2112
$ws.Range("M80").Value = -1114
$ws.Range("H83").Value = 2811.3333
$ws.Range("I83").Value = 2112
$ws.Range("K83").Value = 10560
$ws.Range("M83").Value = -5568

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1674
$ws.Range("I7").Value = 1556.1428
$ws.Range("K7").Value = 1556.1428
$ws.Range("M7").Value = -1444.1428
$ws.Range("H40").Value = 925.8570999999999
$ws.Range("I40").Value = 896.2
$ws.Range("K40").Value = 896.2
$ws.Range("M40").Value = -760.2
$ws.Range("H55").Value = 439.7857
$ws.Range("I55").Value = 308.1579
$ws.Range("K55").Value = 308.1579
$ws.Range("M55").Value = -135.1579
$ws.Range("H56").Value = 13025.5
$ws.Range("I56").Value = 16051
$ws.Range("J56").Value = 10000
$ws.Range("K56").Value = 16051
$ws.Range("L56").Value = 10000
$ws.Range("M56").Value = -15360
$ws.Range("N56").Value = -11382
$ws.Range("H68").Value = 2354
$ws.Range("I68").Value = 2284.8
$ws.Range("K68").Value = 2284.8
$ws.Range("M68").Value = -1535.8
$ws.Range("H71").Value = 2354
$ws.Range("I71").Value = 2284.8
$ws.Range("K71").Value = 11424
$ws.Range("M71").Value = -7680
$ws.Range("H93").Value = 303
$ws.Range("I93").Value = 303
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 303
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 945
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 4050.25
$ws.Range("I100").Value = 1700.375
$ws.Range("J100").Value = 8750
$ws.Range("K100").Value = 1700.375
$ws.Range("L100").Value = 8750
$ws.Range("M100").Value = -1159.375
$ws.Range("N100").Value = -9832
$ws.Range("H126").Value = 1674
$ws.Range("I126").Value = 1556.1428
$ws.Range("K126").Value = 4668.428400000001
$ws.Range("M126").Value = -2198.428400000001
$ws.Range("H136").Value = 2259.4285
$ws.Range("I136").Value = 2259.4285
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6778.2855
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4228.2855
$ws.Range("N136").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1267.1111
$ws.Range("I122").Value = 771.2857
$ws.Range("K122").Value = 2313.8571
$ws.Range("M122").Value = 136.1428999999998
$ws.Range("H126").Value = 145855
$ws.Range("I126").Value = 200797.2
$ws.Range("K126").Value = 602391.6000000001
$ws.Range("M126").Value = -599921.6000000001

Write-Host "Seraph_Profits: scheduled runner update applied ($($wb.Worksheets.Count) sheets touched: 8)"
